# Update the "想去人数" (number of people interested) figures that changed
# between the two scrapes, on both the "展览" sheet and the combined
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 25
$ws1.Range("F5").Value = 4515
$ws1.Range("F6").Value = 166
$ws1.Range("F8").Value = 264

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 25
$ws4.Range("F9").Value = 4515
$ws4.Range("F10").Value = 166
$ws4.Range("F13").Value = 264
